# Ajout d'une colonne "Ville" au tableau des aéroports.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$tbl = $ws.ListObjects.Item(1)

# Agrandit le tableau d'une colonne (F7 -> G7) en ajoutant une ListColumn.
$newCol = $tbl.ListColumns.Add()

# En-tête + valeurs de la nouvelle colonne "Ville".
$ws.Range("G1").Value = "Ville"
$ws.Range("G2").Value = "Lyon"
$ws.Range("G3").Value = "Longvic"
$ws.Range("G4").Value = "Nancy"
$ws.Range("G5").Value = "Valence"
$ws.Range("G6").Value = "La Rochelle"
$ws.Range("G7").Value = "Niort"

# Mise en page (la sauvegarde a "touché" la mise en page / impression).
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1

# Position du curseur au moment de l'enregistrement.
$ws.Range("H1").Select()
